$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 and A3 with the combined tuple-style strings
$ws.Range("A2").Value = "('Counterspell', ['{U}{U}', 'Instant', 'Counter target spell.'])"
$ws.Range("A3").Value = "('Incinerate', ['{1}{R}', 'Instant', 'Incinerate deals 3 damage to any target. A creature dealt damage this way can’t be regenerated this turn.'])"

# Delete the now-unused rows 4-9 (shift cells up) since the data moved into A2/A3
$ws.Range("A4:A9").Delete(-4162)
